$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 currently has "Confirmed" in the Status column (G9); it should be "Pending"
# to match row 10. This also drops the now-unused "Confirmed" shared string.
$ws.Range("G9").Value = "Pending"
